$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("egresados")

# Fix career (carreras) value for the "Eynar Josue Lanuza" row from
# "Ingeniería de sistemas" to "Ingeniería civil"
$ws.Range("K4").Value = "Ingeniería civil"

# Fix erroneous text entered in the fecha_nacimiento (date) column for
# the "Jorshua Alberto Arauz Cantarero" row
$ws.Range("F3").Value = "asdadsadsasd"

# Update the active selection shown when the workbook was last saved
$ws.Range("F10").Select()
